$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 7, shifting the existing weekly entries (old rows
# 7-27) down to rows 8-28. All of row 27's data ends up duplicated into the
# new row 28, matching the template pattern used for every entry.
$ws.Rows.Item(7).Insert()

# Carry over the date number format used by the rest of column D onto the
# freshly inserted D7 cell before writing its value.
$ws.Range("D7").NumberFormat = $ws.Range("D8").NumberFormat()

# Populate the new weekly price entry in row 7.
$ws.Range("A7").Value = 3
$ws.Range("B7").Value = 'Femacal de La Calera'
$ws.Range("C7").Value = 'Coquimbo'
$ws.Range("D7").Value = 44804
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = 100112044
$ws.Range("G7").Value = 'Perejil'
$ws.Range("H7").Value = 'Sin especificar'
$ws.Range("I7").Value = 'Primera'
$ws.Range("J7").Value = 85
$ws.Range("K7").Value = 3000
$ws.Range("L7").Value = 3000
$ws.Range("M7").Value = 3000
$ws.Range("N7").Value = '$/docena de atados (3 kilos)'
$ws.Range("O7").Value = 'Provincia de Quillota'
$ws.Range("P7").Value = 1000
$ws.Range("Q7").Value = 3
$ws.Range("R7").Value = 'Hortaliza'
